$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.569.48"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "3.112.31"
$ws.Range("E3").Value = "  +1.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.48"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.03"
$ws.Range("E6").Value = "  -1.83%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.097.31"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +1.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("E10").Value = "  +2.31%  "

$ws.Range("E11").Value = "  +5.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.51"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000219"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").Value = "3.603.42"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").Value = "63.440.75"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.112"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").Value = "3.101.14"
$ws.Range("E18").Value = "  +0.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "510.36"
$ws.Range("E19").Value = "  +4.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.75"
$ws.Range("E20").Value = "  +2.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("E22").Value = "  +4.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.51"
$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.31"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  +2.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.30"
$ws.Range("E28").Value = "  -0.71%  "

$ws.Range("E29").Value = "  -0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.48"
$ws.Range("E31").Value = "  +3.04%  "

$ws.Range("E32").Value = "  -3.47%  "

$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "539.22"
$ws.Range("E34").Value = "  -8.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.01"
$ws.Range("E35").Value = "  +13.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.21"
$ws.Range("E37").Value = "  -2.87%  "

$ws.Range("E38").Value = "  +4.48%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0804"
$ws.Range("E39").Value = "  +2.15%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.098.97"
$ws.Range("E40").Value = "  +3.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.17"
$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.66"
$ws.Range("E43").Value = "  -6.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").Value = "  +6.09%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.14"
$ws.Range("E45").Value = "  +2.64%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.54"
$ws.Range("E47").Value = "  +1.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.21"
$ws.Range("E48").Value = "  -3.11%  "

$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("D50").Value = "0.0₃0504"
$ws.Range("E50").Value = "  -5.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.38"
$ws.Range("E51").Value = "  +69.62%  "
